# Generate Report for Archive
# Update localization status for the two files that moved from
# "Ready for handoff" to "In Translation": 04a557fc-... and 8905827f-...
# (the daf95ae8-... row stays "Ready for handoff").

$wb = $excel.ActiveWorkbook

# Overview sheet: columns E (zh-cn) and F (de-de) show per-locale status.
# Row 3 -> 04a557fc-7306-410e-9604-78d8bdb77a69.md
# Row 4 -> 8905827f-75f3-4725-98c1-1e1795ff5034.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"

# zh-cn sheet: column C is "Status".
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C4").Value = "In Translation"

# de-de sheet: column C is "Status".
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C4").Value = "In Translation"
